$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# M1/M2: new header + value (added first so new shared strings land in this order)
$ws.Range("M1").Value = "请求头imprint"
$ws.Range("M2").Value = "oWRkU0X0y2TYFRDqFFdGW153oLM0"

# K2: "杨锋" -> "wStar"
$ws.Range("K2").Value = "wStar"

# Column widths: L (12) widened, M (13) new
# (values chosen so the COM -> stored-width rounding lands on the closest
# achievable match to the target XML widths of 17.5546875 / 38.88671875)
$ws.Columns.Item(12).ColumnWidth = 16.857142857142854
$ws.Columns.Item(13).ColumnWidth = 38.14285714285714

# Selection moves from E15 to D15
[void]$ws.Range("D15").Select()
